# "completing reports for laterality"
# Fill in the pression-diff / mask-pressure / airflow columns (G:I) for the
# rows that were still missing them (data rows 2-58). G and H get a literal
# 0, I gets a formula that mirrors the "mema" column (F) for that row -
# matching the pattern already used for the rows below (59-75) that were
# completed earlier.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow  = 58

# G2:H58 are all literal zeros - set the whole block in one shot.
$ws.Range("G$firstRow`:H$lastRow").Value = 0

# I2:I58 mirror column F on the same row.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("I$r").Formula = "=F$r"
}

# Row 75 had an explicit row height (15.5) left over from formatting; let it
# go back to the sheet's natural/auto height.
$ws.Rows.Item(75).AutoFit()

# Restore the view: no frozen/scrolled top-left cell, selection on H3.
$ws.Activate()
$ws.Range("H3").Select()
